$wb = $excel.ActiveWorkbook

# Insert the new "Products" worksheet right after "LoginPage"
$loginSheet = $wb.Worksheets.Item("LoginPage")
$ws = $wb.Worksheets.Add($null, $loginSheet)
$ws.Name = "Products"

# Populate the new sheet with category data
$ws.Cells.Item(1,1).Value = "panelvalues"
$ws.Cells.Item(1,2).Value = "Variations;Import Products;Import Opening Stock;Selling Price Group;Units;Categories ;Brands"
$ws.Cells.Item(2,1).Value = "categoryName"
$ws.Cells.Item(2,2).Value = "Tea powder"
$ws.Cells.Item(3,1).Value = "categoryCode"
$ws.Cells.Item(3,2).Value = 22

# Formatting to match the source sheet layout (column widths quantize to the
# engine's character-width grid; these inputs land on the 13 / 78.5 cells)
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666
$ws.Columns.Item(2).ColumnWidth = 77.675

$ws.Cells.Item(3,2).HorizontalAlignment = -4131
$ws.Cells.Item(3,2).VerticalAlignment = -4108

# Match the selection/active cell state captured in the workbook
$ws.Range("A3").Select()

Write-Host "Products sheet added"
